$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 511
$ws.Range("F6").Value = 977
$ws.Range("F7").Value = 492
$ws.Range("F9").Value = 2276
$ws.Range("F11").Value = 321
$ws.Range("F12").Value = 133
$ws.Range("F13").Value = 1182
$ws.Range("F15").Value = 2297
$ws.Range("F16").Value = 749
$ws.Range("F17").Value = 17403
$ws.Range("F18").Value = 20
$ws.Range("F19").Value = 1518
$ws.Range("F20").Value = 600
$ws.Range("F22").Value = 267
$ws.Range("F23").Value = 557
$ws.Range("F24").Value = 163
$ws.Range("F25").Value = 135
$ws.Range("F27").Value = 291
$ws.Range("F31").Value = 51

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 158
$ws.Range("F10").Value = 95
$ws.Range("F11").Value = 78
$ws.Range("F12").Value = 49

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5763
$ws.Range("F3").Value = 513
$ws.Range("F4").Value = 502

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 513
$ws.Range("F4").Value = 502
$ws.Range("F5").Value = 977
$ws.Range("F7").Value = 492
$ws.Range("F9").Value = 2276
$ws.Range("F11").Value = 321
$ws.Range("F13").Value = 133
$ws.Range("F15").Value = 1182
$ws.Range("F18").Value = 158
$ws.Range("F20").Value = 2297
$ws.Range("F21").Value = 749
$ws.Range("F22").Value = 17406
$ws.Range("F23").Value = 20
$ws.Range("F24").Value = 95
$ws.Range("F25").Value = 78
$ws.Range("F26").Value = 1518
$ws.Range("F27").Value = 600
$ws.Range("F29").Value = 267
$ws.Range("F30").Value = 557
$ws.Range("F31").Value = 163
$ws.Range("F32").Value = 135
$ws.Range("F37").Value = 291
$ws.Range("F49").Value = 51
